$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STX")

# Row 4: Inventory
$ws.Range("B4").Value = 1281000000.0
$ws.Range("C4").Value = 1318000000.0
$ws.Range("D4").Value = 1323000000.0
$ws.Range("E4").Value = 1142000000.0
$ws.Range("F4").Value = 1102000000.0

# Row 14: Accounts Payable
$ws.Range("B14").Value = 1861000000.0
$ws.Range("C14").Value = 1730000000.0
$ws.Range("D14").Value = 1795000000.0
$ws.Range("E14").Value = 1808000000.0
$ws.Range("F14").Value = 1830000000.0

# Row 21: Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = -1118000000.0
$ws.Range("C21").Value = -1120000000.0
$ws.Range("D21").Value = -1125000000.0
$ws.Range("E21").Value = -1120000000.0
$ws.Range("F21").Value = -1112000000.0

# Row 32: Net Debt
$ws.Range("G32").Value = 2397000000.0

# Row 33: Total Debt
$ws.Range("G33").Value = 4141000000.0

$wb.Save()
